$d = $word.ActiveDocument
$full = $d.Content
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Of all Kickstarter campaigns, between 2009 and 2017, 5</w:t></w:r><w:r><w:t>3</w:t></w:r><w:r><w:t xml:space="preserve">% are successfully funded. </w:t></w:r><w:r><w:t xml:space="preserve">Theater projects make up 34% of campaigns, successful or otherwise. 77% of the theater campaigns are plays, of which 59% are based in the United States. No journalism projects have been successfully funded during this </w:t></w:r><w:r><w:t>period</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>The</w:t></w:r><w:r><w:t>re is no information on how much each backer pledged to a campaign, as well as no explanation of the difference between a &#8220;failed&#8221; campaign and a &#8220;canceled&#8221; one. The addition of &#8220;live&#8221; campaigns means that there are incomplete data mixed into the data set, which has the potential to skew the above conclusions.</w:t></w:r></w:p><w:p><w:r><w:t>Other tables/graphs which could be generated are  &#8220;&#8216;Staff Picked&#8217; campaigns which were</w:t></w:r><w:r><w:t xml:space="preserve"> successful, failed, canceled, or are currently live</w:t></w:r><w:r><w:t>&#8221; and the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>&#8220;</w:t></w:r><w:r><w:t>&#8216;</w:t></w:r><w:r><w:t>Spotlight</w:t></w:r><w:r><w:t>&#8217;</w:t></w:r><w:r><w:t xml:space="preserve"> campaigns which were successful, failed, canceled, or are currently live</w:t></w:r><w:r><w:t>&#8221;</w:t></w:r><w:r><w:t>, as those could show correlations between the state of each campaign within those particular campaign identifiers.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>BONUS</w:t></w:r></w:p><w:p><w:r><w:t>The median summarizes the campaign data more than the mean, as when considering both sets of data, the &#8220;successful&#8221; campaigns far outweighs the &#8220;failed&#8221; ones.</w:t></w:r></w:p><w:p><w:r><w:t>There is more variance with successful campaigns</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>as the data shows the wide range between the minimum and maximum numbers of backers.</w:t></w:r><w:r><w:t xml:space="preserve"> This makes sense because the more backers a campaign has, the more likely it is to be successful.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full.InsertXML($xml)
Write-Output ("ParaCount=" + $d.Paragraphs.Count)
